$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.069.76'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '2.601.65'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '312.49'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').Value = '98.83'
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('D7').Value = '0.599'
$ws.Range('E7').Value = '  -1.03%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.582'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').Value = '39.03'
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('D11').Value = '54.51'
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('D12').Value = '0.0841'
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Value = '8.14'
$ws.Range('E13').Value = '  -1.71%  '
$ws.Range('D14').Value = '2.993.81'
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('D16').Value = '2.591.24'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').Value = '0.916'
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').Value = '14.88'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').Value = '46.168.66'
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '6.73'
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').Value = '12.81'
$ws.Range('E22').Value = '  -4.17%  '
$ws.Range('D23').Value = '290.39'
$ws.Range('E23').Value = '  +12.93%  '
$ws.Range('D24').Value = '72.65'
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('D25').Value = '3.07'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').Value = '2.27'
$ws.Range('E26').Value = '  +1.80%  '
$ws.Range('D27').Value = '30.11'
$ws.Range('E27').Value = '  +6.43%  '
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('D30').Value = '10.79'
$ws.Range('E30').Value = '  +2.29%  '
$ws.Range('D31').Value = '2.21'
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('D32').Value = '38.14'
$ws.Range('E32').Value = '  -4.11%  '
$ws.Range('D33').Value = '6.26'
$ws.Range('E33').Value = '  +1.00%  '
$ws.Range('D34').Value = '3.59'
$ws.Range('E34').Value = '  -3.81%  '
$ws.Range('D35').Value = '155.75'
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('D36').Value = '0.0840'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '2.21'
$ws.Range('E37').Value = '  -5.73%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').Value = '2.81'
$ws.Range('E38').Value = '  -4.58%  '
$ws.Range('D39').Value = '0.122'
$ws.Range('E39').Value = '  +3.21%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = '15.83'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '22.27'
$ws.Range('E42').Value = '  +16.19%  '
$ws.Range('D43').Value = '0.0331'
$ws.Range('E43').Value = '  +1.56%  '
$ws.Range('D44').Value = '3.59'
$ws.Range('E44').Value = '  -1.72%  '
$ws.Range('E45').Value = '  -6.19%  '
$ws.Range('D46').Value = '2.098.22'
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('D47').Value = '97.69'
$ws.Range('E47').Value = '  +6.74%  '
$ws.Range('D48').Value = '0.998'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').Value = '9.65'
$ws.Range('E49').Value = '  +4.77%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.202'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '108.54'
$ws.Range('E51').Value = '  -2.00%  '
